# Auto-generated PowerShell Excel COM-interop script
# Applies 110 numeric cell updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "H2" = 3.4
    "H3" = 2.38
    "F4" = 2.86
    "O4" = 1.31
    "H5" = 1.19
    "P6" = 1.52
    "F7" = 3.35
    "G7" = 5.3
    "H7" = 1.92
    "K7" = 5.7
    "I12" = 3.85
    "P12" = 1.78
    "Q12" = 2.02
    "G15" = 2.88
    "I15" = 3.1
    "J15" = 3.4
    "N15" = 3.6
    "Q15" = 1.78
    "W15" = 1.53
    "AB15" = 15.5
    "AE15" = 36
    "T18" = 1.74
    "X19" = 21
    "AI19" = 980
    "AM19" = 70
    "H20" = 3.25
    "F21" = 1.68
    "G21" = 1.99
    "H21" = 3.7
    "K21" = 7.6
    "W21" = 2
    "H22" = 2.02
    "K22" = 4
    "R22" = 1.59
    "V22" = 1.86
    "R24" = 1.62
    "S24" = 2.4
    "F25" = 3.1
    "G25" = 4.3
    "H27" = 2.46
    "T27" = 1.89
    "U27" = 1.89
    "V27" = 1.57
    "I29" = 2.4
    "Q29" = 1.79
    "T29" = 1.66
    "U29" = 2.2
    "AK29" = 46
    "H30" = 3.25
    "K30" = 950
    "N32" = 2.34
    "Q33" = 1.92
    "S33" = 3.4
    "R35" = 1.52
    "S35" = 2.18
    "S37" = 3.45
    "F39" = 3.25
    "AC39" = 8.800000000000001
    "AK39" = 1000
    "N42" = 2.98
    "H43" = 1.78
    "I43" = 1.91
    "Q43" = 1.45
    "R43" = 1.62
    "Y43" = 16
    "H44" = 3.65
    "S44" = 1.92
    "Z45" = 32
    "AE45" = 40
    "AJ45" = 36
    "S47" = 2.52
    "U47" = 2.46
    "V47" = 1.4
    "N48" = 5.5
    "F49" = 3.3
    "G49" = 4.7
    "H49" = 1.9
    "I49" = 2.2
    "K49" = 5.9
    "V49" = 1.83
    "G51" = 5.7
    "H51" = 1.74
    "I51" = 1.91
    "J51" = 3.95
    "K51" = 5.4
    "W51" = 1.25
    "AO51" = 13
    "F52" = 1.77
    "G52" = 1.82
    "K52" = 3.9
    "V52" = 1.2
    "W52" = 2.2
    "AC52" = 9.800000000000001
    "AF52" = 10
    "AN52" = 13
    "F56" = 1.75
    "G56" = 1.8
    "H56" = 5.5
    "S56" = 3.65
    "W56" = 2.24
    "N58" = 3.65
    "G59" = 3.6
    "I59" = 3.45
    "V59" = 1.43
    "W59" = 1.34
    "F62" = 3.7
    "I62" = 2.14
    "V62" = 1.89
    "Q64" = 1.52
    "R64" = 1.43
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
